$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.948.31"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.121.75"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'578.38"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'172.16"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'0.0000248"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'37.17"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "3.638.47"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "66.904.76"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "3.121.31"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'16.25"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'475.36"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").Value = "'0.709"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("D23").Value = "'83.77"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'13.26"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "'2.38"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "'0.115"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "0.0₃0948"
$ws.Range("E33").Value = "  -6.78%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("D37").Value = "'47.03"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "2.817.99"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "'383.38"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "'0.0353"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("E46").Value = "  -9.78%  "
$ws.Range("D47").Value = "'135.68"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D49").Value = "'24.95"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("E51").Value = "  -0.85%  "
